# Update cryptos list: refreshed prices / 1h volume percentages,
# plus a Chainlink / Solana / BinanceUSD row reshuffle (rows 13-15).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B-E on the data rows are stored as text (prices like "28.530.60"
# or "0.00001088" are NOT numbers) -- force text format first so COM
# does not coerce the assigned string into a Double and mangle it
# (dropping trailing zeros, switching to scientific notation, etc).
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '28.530.60'
$ws.Range('E2').Value = '  +1.50%  '

$ws.Range('D3').Value = '1.827.00'
$ws.Range('E3').Value = '  +1.81%  '

$ws.Range('E4').Value = '  +0.00%  '

$ws.Range('D5').Value = '317.34'
$ws.Range('E5').Value = '  +0.00%  '

$ws.Range('D6').Value = '1.001'
$ws.Range('E6').Value = '  +0.06%  '

$ws.Range('D7').Value = '0.5436'
$ws.Range('E7').Value = '  +0.32%  '

$ws.Range('D8').Value = '0.4036'
$ws.Range('E8').Value = '  +6.75%  '

$ws.Range('D9').Value = '0.07679'
$ws.Range('E9').Value = '  +3.06%  '

$ws.Range('E10').Value = '  +2.28%  '

$ws.Range('D11').Value = '41.88'
$ws.Range('E11').Value = '  +0.36%  '

$ws.Range('D12').Value = '6.329'
$ws.Range('E12').Value = '  +3.45%  '

$ws.Range('B13').Value = 'Solana'
$ws.Range('C13').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D13').Value = '21.09'
$ws.Range('E13').Value = '  +2.62%  '

$ws.Range('B14').Value = 'BinanceUSD'
$ws.Range('C14').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D14').Value = '1.001'
$ws.Range('E14').Value = '  -0.01%  '

$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D15').Value = '7.613'
$ws.Range('E15').Value = '  +5.05%  '

$ws.Range('D16').Value = '1.824.07'
$ws.Range('E16').Value = '  +1.98%  '

$ws.Range('D17').Value = '0.00001088'
$ws.Range('E17').Value = '  +2.74%  '

$ws.Range('D18').Value = '89.97'
$ws.Range('E18').Value = '  +0.87%  '

$ws.Range('D19').Value = '0.06607'
$ws.Range('E19').Value = '  +1.87%  '

$ws.Range('D20').Value = '17.81'
$ws.Range('E20').Value = '  +2.97%  '

$ws.Range('E21').Value = '  +0.02%  '

$ws.Range('D22').Value = '6.075'
$ws.Range('E22').Value = '  +2.89%  '

$ws.Range('D23').Value = '28.534.37'
$ws.Range('E23').Value = '  +1.45%  '

$ws.Range('D24').Value = '11.19'
$ws.Range('E24').Value = '  +0.09%  '

$ws.Range('D25').Value = '2.279'
$ws.Range('E25').Value = '  +9.02%  '

$ws.Range('D26').Value = '157.99'
$ws.Range('E26').Value = '  +2.02%  '

$ws.Range('D27').Value = '20.79'
$ws.Range('E27').Value = '  +2.55%  '

$ws.Range('D28').Value = '2.451'
$ws.Range('E28').Value = '  +6.87%  '

$ws.Range('D29').Value = '2.035.25'
$ws.Range('E29').Value = '  +2.06%  '

$ws.Range('D30').Value = '124.03'
$ws.Range('E30').Value = '  +2.47%  '

$ws.Range('D31').Value = '1.125'
$ws.Range('E31').Value = '  +0.38%  '

$ws.Range('D32').Value = '0.1107'
$ws.Range('E32').Value = '  +4.76%  '

$ws.Range('D33').Value = '5.680'
$ws.Range('E33').Value = '  +2.18%  '

$ws.Range('D34').Value = '0.07373'
$ws.Range('E34').Value = '  +13.16%  '

$ws.Range('D35').Value = '3.644'
$ws.Range('E35').Value = '  -0.22%  '

$ws.Range('D36').Value = '0.2241'
$ws.Range('E36').Value = '  -0.85%  '

$ws.Range('D37').Value = '0.02355'
$ws.Range('E37').Value = '  +2.65%  '

$ws.Range('D38').Value = '5.214'
$ws.Range('E38').Value = '  +3.76%  '

$ws.Range('D39').Value = '8.901'
$ws.Range('E39').Value = '  +5.21%  '

$ws.Range('D40').Value = '0.6300'
$ws.Range('E40').Value = '  +1.99%  '

$ws.Range('E41').Value = '  +2.45%  '

$ws.Range('D42').Value = '1.190'

$ws.Range('D43').Value = '1.000'
$ws.Range('E43').Value = '  +0.03%  '

$ws.Range('D44').Value = '1.401'
$ws.Range('E44').Value = '  -3.52%  '

$ws.Range('D45').Value = '13.43'
$ws.Range('E45').Value = '  +0.42%  '

$ws.Range('D46').Value = '0.5879'
$ws.Range('E46').Value = '  +1.53%  '

$ws.Range('D47').Value = '3.706'
$ws.Range('E47').Value = '  +0.88%  '

$ws.Range('D48').Value = '125.38'
$ws.Range('E48').Value = '  +0.17%  '

$ws.Range('D49').Value = '2.003'
$ws.Range('E49').Value = '  +4.20%  '

$ws.Range('D50').Value = '1.198'
$ws.Range('E50').Value = '  +0.65%  '

$ws.Range('D51').Value = '0.06907'
$ws.Range('E51').Value = '  +1.42%  '

Write-Output "cryptos updated"